$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, shifting existing rows 122:183 down to 123:184.
$ws.Rows.Item(122).Insert()

# Populate the new row 122 with the new data record.
$ws.Range("A122").Value = 9
$ws.Range("B122").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C122").Value = "Metropolitana"
$ws.Range("D122").Value = 44529
$ws.Range("D122").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E122").Value = 13
$ws.Range("F122").Value = 100112026
$ws.Range("G122").Value = "Haba"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 34
$ws.Range("K122").Value = 7000
$ws.Range("L122").Value = 8000
$ws.Range("M122").Value = 7500
$ws.Range("N122").Value = "$/saco 25 kilos"
$ws.Range("O122").Value = "Carahue"
$ws.Range("P122").Value = 300
$ws.Range("Q122").Value = 25
$ws.Range("R122").Value = "Hortaliza"
